# Update "想去人数" (interest count) figures in column F for the
# "展览" (Exhibition) sheet and the combined "全部类型" (All Types) sheet.
#
# Row -> (old, new) changes are identical in content between the two
# sheets, except that "全部类型" has one extra row (a 演出 / Performance
# entry) inserted at row 22, shifting everything below it down by one.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F3").Value  = 218
$wsExpo.Range("F5").Value  = 6618
$wsExpo.Range("F8").Value  = 133
$wsExpo.Range("F9").Value  = 6010
$wsExpo.Range("F10").Value = 41
$wsExpo.Range("F11").Value = 191
$wsExpo.Range("F14").Value = 86
$wsExpo.Range("F15").Value = 386
$wsExpo.Range("F16").Value = 93
$wsExpo.Range("F17").Value = 17
$wsExpo.Range("F18").Value = 351
$wsExpo.Range("F19").Value = 41
$wsExpo.Range("F21").Value = 4330
$wsExpo.Range("F22").Value = 45
$wsExpo.Range("F23").Value = 9
$wsExpo.Range("F24").Value = 187
$wsExpo.Range("F25").Value = 26

# --- Sheet "全部类型" ---------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F3").Value  = 218
$wsAll.Range("F5").Value  = 6618
$wsAll.Range("F8").Value  = 133
$wsAll.Range("F9").Value  = 6010
$wsAll.Range("F10").Value = 41
$wsAll.Range("F11").Value = 191
$wsAll.Range("F14").Value = 86
$wsAll.Range("F15").Value = 386
$wsAll.Range("F16").Value = 93
$wsAll.Range("F17").Value = 17
$wsAll.Range("F18").Value = 351
$wsAll.Range("F19").Value = 41
$wsAll.Range("F21").Value = 4330
$wsAll.Range("F23").Value = 45
$wsAll.Range("F24").Value = 9
$wsAll.Range("F25").Value = 187
$wsAll.Range("F26").Value = 26
